$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("codigos")
$ws2 = $wb.Worksheets.Item("tipos infraccion")

# --- sheet2 "tipos infraccion": shift the existing "Tipo" rows (3..17) down
# one row to make room for the new "Alcoholemia y otras Negativo" row, then
# fill in the new abbreviation column (B) and the new D2 value.
# Shift bottom-up so we never clobber a row before it's been read.
for ($r = 17; $r -ge 3; $r--) {
    $dest = $r + 1
    $txt = $ws2.Cells.Item($r, 1).Text
    $ws2.Cells.Item($dest, 1).Value = $txt
}

# New abbreviation header + values for the (still pre-shift-positioned)
# rows 1 and 2.
$ws2.Range("B1").Value = "Abreviacion"
$ws2.Range("B2").Value = "ALCO"

# Abbreviations for the rows that were just shifted down to 4..18.
$ws2.Range("B4").Value = "BARR"
$ws2.Range("B5").Value = "CARR"
$ws2.Range("B6").Value = "CELU"
$ws2.Range("B7").Value = "COND"
$ws2.Range("B8").Value = "DOCU"
$ws2.Range("B9").Value = "ENSE"
$ws2.Range("B10").Value = "ESTA"
$ws2.Range("B11").Value = "INDI"
$ws2.Range("B12").Value = "LICE"
$ws2.Range("B13").Value = "LUCE"
$ws2.Range("B14").Value = "LUZ "
$ws2.Range("B15").Value = "MOTO"
$ws2.Range("B16").Value = "PEAT"
$ws2.Range("B17").Value = "PRIN"
$ws2.Range("B18").Value = "SEGU"

# Extra note in D2.
$ws2.Range("D2").Value = "ALCO-DOC-LICE"

# New row 3: "Alcoholemia y otras Negativo" / "ALCN".
$ws2.Range("A3").Value = "Alcoholemia y otras Negativo"
$ws2.Range("B3").Value = "ALCN"

# --- view/selection bookkeeping to match the saved workbook state.
$ws1.Activate() | Out-Null
$ws1.Range("C3").Select() | Out-Null

$ws2.Activate() | Out-Null
$ws2.Range("F8").Select() | Out-Null
